$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 0.769
$ws.Range("G8").Value = 0.746
$ws.Range("H8").Formula = "=F8-G8"

$ws.Range("H8").Select()
